$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C1 (was S100A8) -> CD14, D1 (was GZMA) -> NKG7
$ws.Range("C1").Value = "CD14"
$ws.Range("D1").Value = "NKG7"

# Update data values that changed
$ws.Range("D2").Value = 1
$ws.Range("C3").Value = 18
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 160

# Remove the per-cell style (s="1", applyFill) so cells revert to default style
$ws.Range("A1:E5").Style = "Normal"
